$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(8692768038931, "CASTROL EDGE 5W/40 1LT TU", "Castrol", "MADENİ YAĞLAR", 77, 3, "ad", 315.006, 315.006),
    @(4008177151385, "CASTROL EDGE 5W40 4 LT", "Castrol", "MADENİ YAĞLAR", 513.7, 2, "ad", 973.9403703703703, 973.9403703703703),
    @(4008177132476, "CASTROL EDGE SUPERCAR TU 10W-60 4LT", "Castrol", "MADENİ YAĞLAR", 729.3, 21, "ad", 1018.794, 1018.794),
    @(4008177157134, "CASTROL ENGİNE SHAMPOO DİESEL 300 ML", "Castrol", "Oto Bakım", 38.5, 34, "ad", 117, 117)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 7).Value = $item[5]
    $ws.Cells.Item($row, 8).Value = $item[6]
    $ws.Cells.Item($row, 9).Value = $item[7]
    $ws.Cells.Item($row, 10).Value = $item[8]
    $row++
}
